$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1700
$ws.Range("J32").Value = 1700
$ws.Range("L32").Value = 1700
$ws.Range("N32").Value = -2352

$ws.Range("H33").Value = 559.9231
$ws.Range("I33").Value = 559.9231
$ws.Range("K33").Value = 559.9231
$ws.Range("M33").Value = -330.9231

$ws.Range("H40").Value = 2904.6
$ws.Range("I40").Value = 2830
$ws.Range("J40").Value = 2979.2
$ws.Range("K40").Value = 2830
$ws.Range("L40").Value = 2979.2
$ws.Range("M40").Value = -2655
$ws.Range("N40").Value = -3329.2

$ws.Range("H88").Value = 1691.091
$ws.Range("I88").Value = 900
$ws.Range("J88").Value = 1770.2
$ws.Range("K88").Value = 900
$ws.Range("L88").Value = 1770.2
$ws.Range("M88").Value = -494
$ws.Range("N88").Value = -2582.2

$ws.Range("H91").Value = 1691.091
$ws.Range("I91").Value = 900
$ws.Range("J91").Value = 1770.2
$ws.Range("K91").Value = 900
$ws.Range("L91").Value = 1770.2
$ws.Range("M91").Value = 504
$ws.Range("N91").Value = -4578.2

$ws.Range("H92").Value = 52631936
$ws.Range("I92").Value = 71428904
$ws.Range("J92").Value = 420.2
$ws.Range("K92").Value = 71428904
$ws.Range("L92").Value = 420.2
$ws.Range("M92").Value = -71427656
$ws.Range("N92").Value = -2916.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 499
$ws.Range("I3").Value = 499
$ws.Range("K3").Value = 499
$ws.Range("M3").Value = -384

$ws.Range("H24").Value = 64000
$ws.Range("J24").Value = 64000
$ws.Range("L24").Value = 64000
$ws.Range("N24").Value = -64748

$ws.Range("H35").Value = 1531.8889
$ws.Range("I35").Value = 1660.875
$ws.Range("K35").Value = 1660.875
$ws.Range("M35").Value = -1254.875

$ws.Range("H100").Value = 64000
$ws.Range("J100").Value = 64000
$ws.Range("L100").Value = 64000
$ws.Range("N100").Value = -66164

$ws.Range("H122").Value = 2638.3076
$ws.Range("I122").Value = 2733.1667
$ws.Range("K122").Value = 8199.500100000001
$ws.Range("M122").Value = -5749.500100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 1660.6
$ws.Range("I7").Value = 434.66666
$ws.Range("J7").Value = 3499.5
$ws.Range("K7").Value = 434.66666
$ws.Range("L7").Value = 3499.5
$ws.Range("M7").Value = -321.66666
$ws.Range("N7").Value = -3725.5

$ws.Range("H8").Value = 53.5
$ws.Range("I8").Value = 54
$ws.Range("J8").Value = 53
$ws.Range("K8").Value = 54
$ws.Range("L8").Value = 53
$ws.Range("M8").Value = 86
$ws.Range("N8").Value = -333

$ws.Range("H86").Value = 3924
$ws.Range("I86").Value = 3706.6667
$ws.Range("J86").Value = 4250
$ws.Range("K86").Value = 3706.6667
$ws.Range("L86").Value = 4250
$ws.Range("M86").Value = -2583.6667
$ws.Range("N86").Value = -6496

$ws.Range("H89").Value = 3924
$ws.Range("I89").Value = 3706.6667
$ws.Range("J89").Value = 4250
$ws.Range("K89").Value = 18533.3335
$ws.Range("L89").Value = 21250
$ws.Range("M89").Value = -12917.3335
$ws.Range("N89").Value = -32482

$ws.Range("H92").Value = 138552.14
$ws.Range("J92").Value = 138552.14
$ws.Range("L92").Value = 138552.14
$ws.Range("N92").Value = -143544.14

$ws.Range("H105").Value = 1882.2
$ws.Range("I105").Value = 1725
$ws.Range("K105").Value = 1725
$ws.Range("M105").Value = 22

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()  # was -5283

$ws.Range("H16").Value = 783
$ws.Range("J16").Value = 800
$ws.Range("L16").Value = 800
$ws.Range("N16").Value = -1374

$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()  # was -43952

$ws.Range("H62").Value = 4966.6665
$ws.Range("I62").Value = 4966.6665
$ws.Range("K62").Value = 4966.6665
$ws.Range("M62").Value = -4342.6665

$ws.Range("H65").Value = 4966.6665
$ws.Range("I65").Value = 4966.6665
$ws.Range("K65").Value = 24833.3325
$ws.Range("M65").Value = -21713.3325

$ws.Range("H86").Value = 333335170
$ws.Range("I86").Value = 500001500
$ws.Range("J86").Value = 2500
$ws.Range("K86").Value = 500001500
$ws.Range("L86").Value = 2500
$ws.Range("M86").Value = -500000377
$ws.Range("N86").Value = -4746

$ws.Range("H89").Value = 333335170
$ws.Range("I89").Value = 500001500
$ws.Range("J89").Value = 2500
$ws.Range("K89").Value = 2500007500
$ws.Range("L89").Value = 12500
$ws.Range("M89").Value = -2500001884
$ws.Range("N89").Value = -23732

$ws.Range("H105").Value = 1895
$ws.Range("I105").Value = 1895
$ws.Range("K105").Value = 1895
$ws.Range("M105").Value = -148

$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()  # was -102954

$ws.Range("H113").Value = 783
$ws.Range("J113").Value = 800
$ws.Range("L113").Value = 800
$ws.Range("N113").Value = -5140

$ws.Range("H122").Value = 1051.7778
$ws.Range("I122").Value = 1109.8572
$ws.Range("K122").Value = 3329.5716
$ws.Range("M122").Value = -879.5715999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 836.7714
$ws.Range("I44").Value = 222
$ws.Range("J44").Value = 916.0968
$ws.Range("K44").Value = 666
$ws.Range("L44").Value = 2748.2904
$ws.Range("M44").Value = -268
$ws.Range("N44").Value = -3544.2904

$ws.Range("H75").Value = 1857.5
$ws.Range("I75").Value = 831.3333
$ws.Range("J75").Value = 2297.2856
$ws.Range("K75").Value = 2493.9999
$ws.Range("L75").Value = 6891.8568
$ws.Range("M75").Value = -1495.9999
$ws.Range("N75").Value = -8887.856800000001

$ws.Range("H78").Value = 1857.5
$ws.Range("I78").Value = 831.3333
$ws.Range("J78").Value = 2297.2856
$ws.Range("K78").Value = 7481.9997
$ws.Range("L78").Value = 20675.5704
$ws.Range("M78").Value = -2489.9997
$ws.Range("N78").Value = -30659.5704

$ws.Range("H103").Value = 399
$ws.Range("I103").Value = 500
$ws.Range("J103").Value = 373.75
$ws.Range("K103").Value = 1500
$ws.Range("L103").Value = 1121.25
$ws.Range("M103").Value = -621
$ws.Range("N103").Value = -2879.25

$ws.Range("H114").Value = 1257.0834
$ws.Range("J114").Value = 1155.7
$ws.Range("L114").Value = 3467.1
$ws.Range("N114").Value = -9975.1

$ws.Range("H117").Value = 2526
$ws.Range("I117").Value = 498.6
$ws.Range("J117").Value = 3447.5454
$ws.Range("K117").Value = 1495.8
$ws.Range("L117").Value = 10342.6362
$ws.Range("M117").Value = 1946.2
$ws.Range("N117").Value = -17226.6362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 527.5
$ws.Range("I5").Value = 527.5
$ws.Range("K5").Value = 527.5
$ws.Range("M5").Value = -415.5

$ws.Range("H94").Value = 62000
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 62000
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 62000
$ws.Range("M94").ClearContents()  # was -24324
$ws.Range("N94").Value = -63352

$ws.Range("H102").Value = 1174.5
$ws.Range("J102").Value = 2007
$ws.Range("L102").Value = 2007
$ws.Range("N102").Value = -5251

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 8105.4443
$ws.Range("I68").Value = 5862.25
$ws.Range("K68").Value = 5862.25
$ws.Range("M68").Value = -5113.25

$ws.Range("H71").Value = 8105.4443
$ws.Range("I71").Value = 5862.25
$ws.Range("K71").Value = 29311.25
$ws.Range("M71").Value = -25567.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2
$ws.Range("I2").Value = 2
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = 110
